$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 16347
$ws.Range("D3").Value = 9211
$ws.Range("E3").Value = 454
$ws.Range("D4").Value = 18258
$ws.Range("D5").Value = 9851
$ws.Range("D6").Value = 196245
$ws.Range("E6").Value = 6813
$ws.Range("D7").Value = 20578
$ws.Range("D8").Value = 16231
$ws.Range("D9").Value = 231802
$ws.Range("E9").Value = 5878
$ws.Range("D10").Value = 35621
$ws.Range("D11").Value = 33238
$ws.Range("E11").Value = 1196
$ws.Range("D12").Value = 20589
$ws.Range("D13").Value = 136894
$ws.Range("E13").Value = 4596
$ws.Range("D14").Value = 28392
$ws.Range("D15").Value = 47537
$ws.Range("E15").Value = 1854
$ws.Range("D16").Value = 82474
$ws.Range("E16").Value = 3478
$ws.Range("D17").Value = 24561
$ws.Range("D18").Value = 85657
$ws.Range("E18").Value = 2093
$ws.Range("D19").Value = 24050
$ws.Range("D20").Value = 57141
$ws.Range("E20").Value = 2233
$ws.Range("D21").Value = 28149
$ws.Range("D22").Value = 160124
$ws.Range("E22").Value = 4249
$ws.Range("D23").Value = 51300
$ws.Range("D24").Value = 22521
$ws.Range("D25").Value = 41194
$ws.Range("D26").Value = 25819
$ws.Range("D27").Value = 19260
$ws.Range("D29").Value = 52445
$ws.Range("D30").Value = 13203
$ws.Range("D31").Value = 88679
$ws.Range("E31").Value = 2386
$ws.Range("D32").Value = 116161
$ws.Range("D33").Value = 15768
$ws.Range("D34").Value = 27774
$ws.Range("D35").Value = 96334
$ws.Range("D36").Value = 159080
$ws.Range("E36").Value = 4090
$ws.Range("D37").Value = 114495
$ws.Range("E37").Value = 4109
$ws.Range("D38").Value = 27216
$ws.Range("D39").Value = 131752
$ws.Range("E39").Value = 2962
$ws.Range("D40").Value = 118346
$ws.Range("D41").Value = 551533
$ws.Range("E41").Value = 13084
$ws.Range("D42").Value = 151623
$ws.Range("E42").Value = 5173
$ws.Range("D43").Value = 421787
$ws.Range("E43").Value = 10805
$ws.Range("D44").Value = 310251
$ws.Range("E44").Value = 11869
$ws.Range("D45").Value = 46593
$ws.Range("D46").Value = 205602
$ws.Range("E46").Value = 4504
$ws.Range("D47").Value = 486575
$ws.Range("E47").Value = 16870
$ws.Range("D48").Value = 63435
$ws.Range("E48").Value = 1822
$ws.Range("D49").Value = 1955
$ws.Range("D50").Value = 125962
$ws.Range("E50").Value = 4096
$ws.Range("D51").Value = 58878
$ws.Range("E51").Value = 1914
$ws.Range("D52").Value = 59288
$ws.Range("D53").Value = 60971
$ws.Range("E53").Value = 1031
$ws.Range("D54").Value = 82811
$ws.Range("E54").Value = 2113
$ws.Range("D55").Value = 317633
$ws.Range("E55").Value = 7971
$ws.Range("D56").Value = 8164
$ws.Range("D57").Value = 1061862
$ws.Range("E57").Value = 24652
$ws.Range("D58").Value = 189449
$ws.Range("E58").Value = 7111
$ws.Range("D59").Value = 47485
$ws.Range("D60").Value = 33637
$ws.Range("D61").Value = 93170
$ws.Range("E61").Value = 3253
$ws.Range("D62").Value = 5993
$ws.Range("D63").Value = 25263
$ws.Range("E63").Value = 754
$ws.Range("D64").Value = 29434
$ws.Range("D65").Value = 551302
$ws.Range("E65").Value = 18000
$ws.Range("D66").Value = 8073
$ws.Range("D67").Value = 9457
$ws.Range("D68").Value = 3905
$ws.Range("E68").Value = 135
$ws.Range("D69").Value = 110540
$ws.Range("D70").Value = 2492
$ws.Range("D72").Value = 7836
$ws.Range("D73").Value = 3580
$ws.Range("D74").Value = 2240
$ws.Range("D75").Value = 51513
$ws.Range("D76").Value = 3220
$ws.Range("D79").Value = 4598
$ws.Range("D80").Value = 387908
$ws.Range("E80").Value = 5373
$ws.Range("D81").Value = 27832
$ws.Range("D83").Value = 731204
$ws.Range("E83").Value = 16134
$ws.Range("D84").Value = 55267
$ws.Range("D85").Value = 160048
$ws.Range("E85").Value = 2995
$ws.Range("F85").Value = 107
